# Auto-generated PowerShell COM-interop script
# Applies the "atualizacao modelos arc e iia" edit:
#  - Appends 20 new Feature/Target rows (421-440) to sheet "Pagina1"
#  - New rows use a Calibri 11 font (Target column right-aligned)
#  - Column A width widened to fit the new, much longer Feature text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$features = @(
  "Quero agendar uma consulta",
  "Eu gostaria de agendar uma consulta para meu filho",
  "Você não atende depois das 17:00 hrs né ?",
  "Qual horário vc tem disponível ? Dia 04/03",
  "Pode ser no dia 07/03",
  "E no dia 12/03 tem horário disponível ?",
  "Ele sai da escola 17:30 teria que ser depois desse horário",
  "O endereço é esse que está no WhatsApp ?",
  "Então eu queria saber o valor da limpeza? Pq aí já agendava a consulta",
  "Infantil",
  "Tá , qual dia tem horário ?",
  "Limpeza , queria para o dia 10 de março",
  "Está confirmado amanhã às 08 horas?",
  "Daqui a pouco estou ai",
  "Bom dia Eni a tia não tá bem tô com muita tosse queria que vc marca se pra outro dia a minha visita aí desculpa beijos 😔",
  "Fui tô tomando remédio obrigada",
  "Estou precisando fazer uma limpeza pesada kkkkkk `n Saiu minha contenção dos dentes de baixo e está até com tártaro 🙆🏻‍♀️",
  "Qual seu último horário de atendimento",
  "Quebrou um dente 😞",
  "Gostaria de agendar um horário"
)

$targets = @(1,1,0,1,1,1,0,0,0,0,1,1,0,0,1,0,1,0,0,1)

$startRow = 421
$n = $features.Length

for ($i = 0; $i -lt $n; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $features[$i]
  $ws.Cells.Item($r, 2).Value = $targets[$i]
}

# New rows get a Calibri 11 font (Feature col regular, Target col right-aligned)
$lastRow = $startRow + $n - 1
$rngA = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($lastRow, 1))
$rngB = $ws.Range($ws.Cells.Item($startRow, 2), $ws.Cells.Item($lastRow, 2))
$rngA.Font.Name = "Calibri"
$rngA.Font.Size = 11
$rngB.Font.Name = "Calibri"
$rngB.Font.Size = 11
$rngB.HorizontalAlignment = -4152

# Widen column A so the longer Portuguese chat text is readable
$ws.Columns.Item(1).ColumnWidth = 157.0

